$d = $word.ActiveDocument

$xmlFragment = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:xml="http://www.w3.org/XML/1998/namespace">
        <w:body><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">Planiranje sprinta inicira Sprint tako što izlaže posao koji treba obaviti za Sprint. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Planiranje sprinta bavi se sledećim temama:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Tema prva: Zašto je ovaj Sprint vredan?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">Vlasnik proizvoda predlaže kako bi proizvod mogao povećati svoju vrednost i korisnost u trenutnom Sprintu. </w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Tim</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve"> zatim sarađuje kako bi defini</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">ao cilj Sprinta koji govori zašto je Sprint vredan </w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>stejkholderima</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>. Cilj sprinta mora biti finaliz</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>ov</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>an pre kraja planiranja sprinta.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Tema dva: Šta se može uraditi u ovom sprintu?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">Kroz razgovor sa vlasnikom proizvoda, programeri biraju stavke iz </w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Product Backlog-a</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve"> koje će uključiti u trenutni Sprint. </w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">Tim </w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>može poboljšati ove stavke tokom ovog procesa, što povećava razumevanje i samopouzdanje.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Odabir koliko se može završiti unutar Sprinta može biti izazovno. Međutim, što više programeri znaju o svojim prošlim performansama, budućim kapacitetima i definiciji gotovog, to će biti sigurniji u svoje prognoze za sprint.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Treća tema: Kako će izabrani posao biti obavljen?</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">Za svaku odabranu stavku </w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>iz Product Backlog-a</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve">, programeri planiraju posao koji je neophodan za kreiranje prirasta koji zadovoljava definiciju gotovog. Ovo se često radi razlaganjem stavki </w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Product Backlog-a</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve"> na manje radne stavke od jednog dana ili manje. Način na koji se to radi je isključivo na nahođenju programera. Niko drugi im ne govori kako da pretvore stavke </w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Product Backlog-a</w:t></w:r><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t xml:space="preserve"> u povećanje vrednosti.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="sr-Latn-RS"/></w:rPr><w:t>Planiranje sprinta je vremenski ograničeno na maksimalno osam sati za jednomesečni sprint. Za kraće sprintove, događaj je obično kraći.</w:t></w:r></w:p></w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$insertionPoint.InsertXML($xmlFragment) | Out-Null

Write-Host "Paragraphs after edit: " $d.Paragraphs.Count
